$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (the "sum" column) onto the new H1 header
# cell so it reuses the existing bold/bordered/centered style instead of
# creating a new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
